$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in this sheet (some
# entries, e.g. "65.562.00", are not valid numbers at all). A plain
# .Value assignment lets Excel auto-convert numeric-looking text (e.g.
# "580.57") into a real number, which would silently change the stored
# type/precision and the serialized text. Force each target cell to Text
# format first so the value is written verbatim as a string, then
# restore the default "Normal" cell style so no stray number format is
# left behind.
$dPriceCells = @(
    'D2',
    'D3',
    'D5',
    'D6',
    'D8',
    'D9',
    'D10',
    'D13',
    'D15',
    'D16',
    'D17',
    'D18',
    'D19',
    'D20',
    'D21',
    'D22',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D30',
    'D33',
    'D34',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($cellRef in $dPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '65.562.00'
$ws.Range('E2').Value = '  +3.10%  '
$ws.Range('D3').Value = '3.456.94'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '580.57'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').Value = '167.92'
$ws.Range('E6').Value = '  +7.12%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.465.64'
$ws.Range('E8').Value = '  +1.85%  '
$ws.Range('D9').Value = '0.562'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = '4.062.33'
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').Value = '27.55'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('D17').Value = '65.535.30'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').Value = '3.458.50'
$ws.Range('E18').Value = '  +2.35%  '
$ws.Range('D19').Value = '6.23'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('D20').Value = '13.80'
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('D21').Value = '382.71'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('D22').Value = '7.92'
$ws.Range('E22').Value = '  +2.09%  '
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').Value = '71.69'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').Value = '0.521'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('D26').Value = '0.0000119'
$ws.Range('E26').Value = '  +2.24%  '
$ws.Range('D27').Value = '9.95'
$ws.Range('E27').Value = '  +3.34%  '
$ws.Range('D28').Value = '0.180'
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').Value = '6.33'
$ws.Range('E30').Value = '  +6.14%  '
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('D33').Value = '23.24'
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('D34').Value = '7.30'
$ws.Range('E34').Value = '  +5.69%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').Value = '160.28'
$ws.Range('D38').Value = '0.900'
$ws.Range('E38').Value = '  +10.43%  '
$ws.Range('D39').Value = '1.86'
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0739'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '26.24'
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.817.53'
$ws.Range('E42').Value = '  +1.38%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '6.62'
$ws.Range('E43').Value = '  +4.19%  '
$ws.Range('D44').Value = '26.66'
$ws.Range('E44').Value = '  +5.69%  '
$ws.Range('D45').Value = '43.04'
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('D46').Value = '4.46'
$ws.Range('E46').Value = '  +1.04%  '
$ws.Range('E47').Value = '  +7.26%  '
$ws.Range('D48').Value = '0.0307'
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('D49').Value = '345.25'
$ws.Range('E49').Value = '  +5.94%  '
$ws.Range('D50').Value = '1.07'
$ws.Range('E50').Value = '  +4.30%  '
$ws.Range('D51').Value = '32.45'
$ws.Range('E51').Value = '  +8.31%  '

foreach ($cellRef in $dPriceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
